$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

# Plain numeric cells
$ws.Cells.Item($row, 1).Value = 112492939
$ws.Cells.Item($row, 2).Value = 90321

# Plain text cells (not numeric-looking, Excel will store as text with no style change)
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
$ws.Cells.Item($row, 4).Value = "NT"

$ws.Cells.Item($row, 5).Value = 2014

$ws.Cells.Item($row, 6).Value = "Koralltaggsvamp"
$ws.Cells.Item($row, 7).Value = "Hericium coralloides"
$ws.Cells.Item($row, 8).Value = "(Scop.:Fr.) Pers."

# I4 holds the text "1" -- numeric-looking, so force text via apostrophe then
# strip the quote-prefix style Excel adds so the cell keeps the default style.
$ws.Cells.Item($row, 9).Value = "'1"
$ws.Cells.Item($row, 9).Style = "Normal"

# J4, K4, N4 are present but empty text cells. A bare "" assignment is treated
# as "no value" by this engine (the cell/row wouldn't be created), so we force
# an empty text entry with a lone apostrophe, then strip the resulting style
# override back to the default.
$ws.Cells.Item($row, 10).Value = "'"
$ws.Cells.Item($row, 10).Style = "Normal"
$ws.Cells.Item($row, 11).Value = "'"
$ws.Cells.Item($row, 11).Style = "Normal"
$ws.Cells.Item($row, 14).Value = "'"
$ws.Cells.Item($row, 14).Style = "Normal"

$ws.Cells.Item($row, 16).Value = "Klåverödsdammen, Sk"

$ws.Cells.Item($row, 17).Value = 386456
$ws.Cells.Item($row, 18).Value = 6211174
$ws.Cells.Item($row, 19).Value = 50

$ws.Cells.Item($row, 20).Value = "Skåne"
$ws.Cells.Item($row, 21).Value = "Svalöv"
$ws.Cells.Item($row, 22).Value = "Skåne"
$ws.Cells.Item($row, 23).Value = "Konga"

# Y4 / AA4 hold dates written as plain text strings ("2023-08-22"), not Excel
# dates -- force text via apostrophe, then strip the quote-prefix style.
$ws.Cells.Item($row, 25).Value = "'2023-08-22"
$ws.Cells.Item($row, 25).Style = "Normal"
$ws.Cells.Item($row, 27).Value = "'2023-08-22"
$ws.Cells.Item($row, 27).Style = "Normal"

$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false

$ws.Cells.Item($row, 32).Value = "'"
$ws.Cells.Item($row, 32).Style = "Normal"

$ws.Cells.Item($row, 33).Value = $false

$ws.Cells.Item($row, 46).Value = "'"
$ws.Cells.Item($row, 46).Style = "Normal"

$ws.Cells.Item($row, 49).Value = "Kenth Sundgren"
$ws.Cells.Item($row, 50).Value = "Kenth Sundgren"

$ws.Cells.Item($row, 51).Value = "'"
$ws.Cells.Item($row, 51).Style = "Normal"
